# Weekly update: insert a new week of "Espárragos" (asparagus) price
# records at Mercado Mayorista Lo Valledor de Santiago, pushing the
# previously-recorded rows down by 3 (173 -> 196 rows total).
#
# The three new records (date serial 45244 = 2023-11-14) are inserted
# right above the existing block that starts at row 174, so everything
# that used to live in rows 174:193 now lives in rows 177:196.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 174:193 down by inserting 3 fresh rows above row 174.
# Excel's native Insert() shifts everything below and carries the
# row/column formatting (e.g. the date style on column D) along.
$ws.Rows("174:176").Insert()

# Columns that are identical across the three new rows.
$ws.Range("A174:A176").Value = 6
$ws.Range("B174:B176").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C174:C176").Value = "Metropolitana"
$ws.Range("D174:D176").Value = 45244
$ws.Range("E174:E176").Value = 13
$ws.Range("F174:F176").Value = 300000000
$ws.Range("G174:G176").Value = "Espárragos"
$ws.Range("H174:H176").Value = "Sin especificar"
$ws.Range("N174:N176").Value = "$/kilo"
$ws.Range("Q174:Q176").Value = 1
$ws.Range("R174:R176").Value = "Hortaliza"

# Row 174 - Banquete
$ws.Range("I174").Value = "Banquete"
$ws.Range("J174").Value = 2800
$ws.Range("K174").Value = 1600
$ws.Range("L174").Value = 1600
$ws.Range("M174").Value = 1600
$ws.Range("O174").Value = "Provincia de Linares"
$ws.Range("P174").Value = 1600

# Row 175 - Primera
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 7000
$ws.Range("K175").Value = 1400
$ws.Range("L175").Value = 1800
$ws.Range("M175").Value = 1640
$ws.Range("O175").Value = "Provincia de Linares"
$ws.Range("P175").Value = 1640

# Row 176 - Segunda
$ws.Range("I176").Value = "Segunda"
$ws.Range("J176").Value = 1800
$ws.Range("K176").Value = 1200
$ws.Range("L176").Value = 1400
$ws.Range("M176").Value = 1311
$ws.Range("O176").Value = "Provincia de Linares"
$ws.Range("P176").Value = 1311
